# Apply attendance updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Invalid (G) and Absent (H) marked
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Absent (H) marked
$ws.Range("H4").Value = 1

# Row 5: Total Attendance Count (D) and Real (E) marked
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Rows 6-18: Absent (H) marked
for ($r = 6; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
